$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.362.04"
$ws.Range("E2").Value = "  +0.53%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.593.50"
$ws.Range("E3").Value = "  +0.74%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.32%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.60"
$ws.Range("E5").Value = "  +0.94%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.16%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  +0.25%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -0.07%  "

# Row 10 - Solana
$ws.Range("E10").Value = "  -0.50%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.08%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.816.80"

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.604.39"
$ws.Range("E13").Value = "  +2.74%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +0.60%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +1.00%  "

# Row 16 - Litecoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.60"
$ws.Range("E16").Value = "  -0.10%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "26.360.15"
$ws.Range("E17").Value = "  +0.50%  "

# Row 18 - ShibaInu
$ws.Range("D18").Value = "0.0₃0732"
$ws.Range("E18").Value = "  -0.68%  "

# Row 19 - Chainlink
$ws.Range("E19").Value = "  +3.71%  "

# Row 20 - BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "212.70"
$ws.Range("E20").Value = "  +2.89%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.25%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +1.10%  "

# Row 23 - Toncoin
$ws.Range("E23").Value = "  -1.32%  "

# Row 24 - Avalanche
$ws.Range("E24").Value = "  +1.94%  "

# Row 25 - Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.81"
$ws.Range("E25").Value = "  +0.11%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  -0.30%  "

# Row 27 - Cosmos
$ws.Range("E27").Value = "  +0.71%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  -0.45%  "

# Row 29 - EthereumClassic
$ws.Range("E29").Value = "  -0.21%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  -0.07%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +1.13%  "

# Row 32 - Filecoin
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.22"
$ws.Range("E32").Value = "  -0.18%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  +1.56%  "

# Row 34 - Maker
$ws.Range("D34").Value = "1.341.96"
$ws.Range("E34").Value = "  +4.26%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  -1.13%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  +0.00%  "

# Row 37 - LidoDAOToken
$ws.Range("E37").Value = "  +0.37%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  +0.25%  "

# Row 39 - was WEMIXToken, now ARBITRUM
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.820"
$ws.Range("E39").Value = "  +0.80%  "

# Row 40 - was ARBITRUM, now WEMIXToken
$ws.Range("B40").Value = "WEMIXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.04"
$ws.Range("E40").Value = "  -19.54%  "

# Row 41 - FraxShare
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.79"
$ws.Range("E41").Value = "  +4.68%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  -0.26%  "

# Row 43 - MXToken
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.13"
$ws.Range("E43").Value = "  +0.07%  "

# Row 44 - TrustWalletToken
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.764"
$ws.Range("E44").Value = "  -0.84%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").Value = "1.729.01"

# Row 46 - Aave
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.69"
$ws.Range("E46").Value = "  -1.41%  "

# Row 47 - Quant
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.80"
$ws.Range("E47").Value = "  -1.15%  "

# Row 48 - BabyDogeCoin
$ws.Range("E48").Value = "  +2.19%  "

# Row 49 - RenderToken
$ws.Range("E49").Value = "  -2.77%  "

# Row 50 - Algorand
$ws.Range("E50").Value = "  -3.01%  "

# Row 51 - Cronos
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0505"
$ws.Range("E51").Value = "  -0.68%  "
